$d = $word.ActiveDocument

# --- Merge paragraphs 5 and 6 ---------------------------------------------
# The original document has two separate paragraphs:
#   P5: "...en carácter de perito informático,{/peritos}{#integrantes !== []}{#integrantes}"
#   P6: "Junto a {nombreYApellido}, ... .-"
# The target collapses these into a single paragraph (P6's <w:pPr>/bookmark
# disappear, its runs are appended to P5's run sequence). Do this by
# deleting the paragraph mark that currently separates them.
$p5 = $d.Paragraphs.Item(5)
$markStart = $p5.Range.End - 1
$markEnd = $p5.Range.End
$d.Range($markStart, $markEnd).Delete()

# --- Change 1 -------------------------------------------------------------
# ", –Prosecretario– en carácter de Jefe de Gabinete de Informática Forense.{#peritos} "
# -> ", –Prosecretario– en carácter de Jefe de Gabinete de Informática Forense y {#peritos} "
$d.Content.Find.Execute(
    "Jefe de Gabinete de Informática Forense.{#peritos} ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Jefe de Gabinete de Informática Forense y {#peritos} ",
    2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# " en carácter de perito informático,{/peritos}{#integrantes !== []}{#integrantes}Junto a"
# -> " en carácter de perito informático.{/peritos}{#integrantes !== []}{#integrantes} Junto a"
$d.Content.Find.Execute(
    "en carácter de perito informático,{/peritos}{#integrantes !== []}{#integrantes}Junto a",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "en carácter de perito informático.{/peritos}{#integrantes !== []}{#integrantes} Junto a",
    2) | Out-Null

# "matrícula Nº {matricula}," -> "matrícula Nº {legajoOMatricula},"
$d.Content.Find.Execute(
    "matrícula Nº {matricula},",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "matrícula Nº {legajoOMatricula},",
    2) | Out-Null

# --- Change 3 -------------------------------------------------------------
# "Un (01) ordenador portátil, ... {#imei===“”} Sin S/N visible{/}{#imei!==“”} con S/N: {imei} {/}"
# -> "... {#serialNumber===“”} Sin S/N visible{/}{#serialNumber!==“”} con S/N: {serialNumber} {/}"
$d.Content.Find.Execute(
    "{#imei===“”} Sin S/N visible{/}{#imei!==“”} con S/N: {imei} {/}, de su interior se extrae:-",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{#serialNumber===“”} Sin S/N visible{/}{#serialNumber!==“”} con S/N: {serialNumber} {/}, de su interior se extrae:-",
    2) | Out-Null

# --- Change 4 -------------------------------------------------------------
# "Un (01) {tipoDeDisco} ... {#imei!==””} con S/N: {imei},{/}{#imei===””}Sin S/N visible,{/}"
# -> "... {#serialNumber!==””} con S/N: {serialNumber},{/}{#serialNumber===””}Sin S/N visible,{/}"
$d.Content.Find.Execute(
    "{#imei!==””} con S/N: {imei},{/}{#imei===””}Sin S/N visible,{/}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{#serialNumber!==””} con S/N: {serialNumber},{/}{#serialNumber===””}Sin S/N visible,{/}",
    2) | Out-Null
